$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "...batch_files/ folder" + "." (two separate runs) become a
# single run "...batch_files/ folder." (text merge, no wording change).
# ---------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("batch_files/ folder.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rng1.Find.Found) {
    throw "Could not find the 'batch_files/ folder.' text"
}
$matchEnd = $rng1.End

# The "/ folder" run is the 8 characters right before the trailing period.
$folderRun = $d.Range($matchEnd - 9, $matchEnd - 1)
if ($folderRun.Text -ne "/ folder") {
    throw "Unexpected text for folder run: '$($folderRun.Text)'"
}
# Appending "." directly onto this run merges it in (identical formatting),
# producing a single run "/ folder.".
$folderRun.InsertAfter(".")

# The original, now-redundant, standalone "." run got pushed one
# character to the right; delete it.
$oldPeriod = $d.Range($matchEnd, $matchEnd + 1)
if ($oldPeriod.Text -ne ".") {
    throw "Unexpected text for old period run: '$($oldPeriod.Text)'"
}
$oldPeriod.Text = ""

# ---------------------------------------------------------------------
# Change 2: "This batch file will start whenever the PC gets restarted."
# becomes "This batch file will start whenever the PC is turned on."
# with "is turned on" and the trailing "." split into their own runs.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("gets restarted", $true, $false, $false, $false, $false, $true, 1, $false, "", 1)
if (-not $rng2.Find.Found) {
    throw "Could not find the 'gets restarted' text"
}
$replaceStart = $rng2.Start
$rng2.Text = "is turned on"

# Toggling a character format on the newly written text and back forces
# it (and the trailing ".") to live in their own runs instead of being
# re-absorbed into the preceding/following run.
$newRun = $d.Range($replaceStart, $replaceStart + 12)
$newRun.Bold = 1
$newRun.Bold = 0
